# Update the cryptocurrency price/volume table (Sheet1) to the latest
# snapshot values, matching the upstream "Updated cryptos list" GitHub
# Actions commit. Column D ("Price") and column E ("Volume(1h)") are
# plain text in this workbook (not numbers), so any Price value that
# would otherwise be auto-recognised by Excel as a pure number is typed
# with a leading apostrophe to force it to stay text, exactly like the
# existing cells in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.125.66"
$ws.Range("E2").Value = "  -4.62%  "
$ws.Range("D3").Value = "3.277.64"
$ws.Range("E3").Value = "  -5.55%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("B5").Value = "Solana"
$ws.Range("C5").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D5").Value = "'186.64"
$ws.Range("E5").Value = "  -3.02%  "
$ws.Range("B6").Value = "BNB"
$ws.Range("C6").Value = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
$ws.Range("D6").Value = "'558.44"
$ws.Range("E6").Value = "  -3.18%  "
$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("E8").Value = "  -2.65%  "
$ws.Range("D9").Value = "3.273.02"
$ws.Range("E9").Value = "  -5.42%  "
$ws.Range("D10").Value = "'0.189"
$ws.Range("E10").Value = "  -7.70%  "
$ws.Range("D11").Value = "'0.588"
$ws.Range("E11").Value = "  -4.71%  "
$ws.Range("D12").Value = "'47.58"
$ws.Range("E12").Value = "  -7.50%  "
$ws.Range("E13").Value = "  -6.14%  "
$ws.Range("D14").Value = "'8.65"
$ws.Range("E14").Value = "  -5.15%  "
$ws.Range("D15").Value = "'635.04"
$ws.Range("E15").Value = "  -1.04%  "
$ws.Range("D16").Value = "3.802.40"
$ws.Range("E16").Value = "  -5.83%  "
$ws.Range("D17").Value = "65.999.00"
$ws.Range("E17").Value = "  -4.50%  "
$ws.Range("D18").Value = "'17.93"
$ws.Range("E18").Value = "  -1.25%  "
$ws.Range("D20").Value = "3.276.21"
$ws.Range("E20").Value = "  -5.52%  "
$ws.Range("D21").Value = "'11.37"
$ws.Range("E21").Value = "  -7.52%  "
$ws.Range("D22").Value = "'0.906"
$ws.Range("E22").Value = "  -3.83%  "
$ws.Range("D23").Value = "'18.49"
$ws.Range("E23").Value = "  +3.81%  "
$ws.Range("D24").Value = "'107.12"
$ws.Range("E24").Value = "  +8.17%  "
$ws.Range("D25").Value = "'4.92"
$ws.Range("E25").Value = "  -6.68%  "
$ws.Range("D26").Value = "'3.97"
$ws.Range("E26").Value = "  -7.02%  "
$ws.Range("D27").Value = "'2.68"
$ws.Range("E27").Value = "  -6.30%  "
$ws.Range("E28").Value = "  -2.74%  "
$ws.Range("D29").Value = "'8.74"
$ws.Range("E29").Value = "  -6.20%  "
$ws.Range("D30").Value = "'30.38"
$ws.Range("E30").Value = "  -6.17%  "
$ws.Range("D31").Value = "'4.07"
$ws.Range("E31").Value = "  -4.21%  "
$ws.Range("D32").Value = "'6.28"
$ws.Range("E32").Value = "  -6.63%  "
$ws.Range("D33").Value = "'11.06"
$ws.Range("E33").Value = "  -4.57%  "
$ws.Range("E34").Value = "  -3.54%  "
$ws.Range("D35").Value = "'532.88"
$ws.Range("E35").Value = "  +2.53%  "
$ws.Range("D36").Value = "'57.54"
$ws.Range("E36").Value = "  -5.44%  "
$ws.Range("D37").Value = "3.732.32"
$ws.Range("E37").Value = "  +0.78%  "
$ws.Range("E38").Value = "  -0.05%  "
$ws.Range("E39").Value = "  -2.66%  "
$ws.Range("D40").Value = "0.0₃0734"
$ws.Range("E40").Value = "  -7.05%  "
$ws.Range("E42").Value = "  -6.01%  "
$ws.Range("D43").Value = "'3.44"
$ws.Range("E43").Value = "  -1.14%  "
$ws.Range("D44").Value = "'32.85"
$ws.Range("E44").Value = "  -4.08%  "
$ws.Range("D45").Value = "'0.339"
$ws.Range("E45").Value = "  -8.56%  "
$ws.Range("E46").Value = "  -1.92%  "
$ws.Range("D47").Value = "'0.0416"
$ws.Range("E47").Value = "  -5.83%  "
$ws.Range("E48").Value = "  -7.21%  "
$ws.Range("E49").Value = "  -3.30%  "
$ws.Range("E50").Value = "  -0.14%  "
$ws.Range("D51").Value = "'1.26"
$ws.Range("E51").Value = "  +2.01%  "
